$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2902457.32
$ws.Range("C7").Value = -34.67462777809543
$ws.Range("D7").Value = 2928
$ws.Range("E7").Value = 2928
$ws.Range("F7").Value = 991.2764071038251
$ws.Range("G7").Value = 5.662897145812851
